$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("earnings_debt")

$ws.Range("D2").Value = 0.121
$ws.Range("E2").Value = 0.07679999999999999
$ws.Range("G2").Value = 0.3078990592802747
$ws.Range("H2").Value = 0.3078990592802747
$ws.Range("I2").Value = 0.2441391667911005
$ws.Range("J2").Value = 0.155083590069419
$ws.Range("K2").Value = 98
$ws.Range("L2").Value = 0.1463341794833508
$ws.Range("M2").Value = 23.25
$ws.Range("N2").Value = 0.01677368155255754
$ws.Range("O2").Value = 0.2372448979591837
$ws.Range("P2").Value = 20.3
$ws.Range("Q2").Value = 0.01464540797922228
$ws.Range("R2").Value = 0.2071428571428572
$ws.Range("S2").Value = 2.949999999999999
$ws.Range("T2").Value = 0.1268817204301075
$ws.Range("U2").Value = 146.7
$ws.Range("V2").Value = 0.1058365197316211
$ws.Range("W2").Value = 0.1528622679769147
$ws.Range("X2").Value = 0.1129633422907675
$ws.Range("Y2").Value = 0.03989892568614717
$ws.Range("Z2").Value = 1.124223602484472
$ws.Range("AA2").Value = 0.1743486323140673
$ws.Range("AB2").Value = 0.108397913298502
$ws.Range("AC2").Value = 0.06595071901556529
$ws.Range("AD2").Value = 128
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 128
$ws.Range("AG2").Value = -18.69999999999999
$ws.Range("AH2").Value = 0.08453866983686679
$ws.Range("AI2").Value = 0.127987201279872
$ws.Range("AJ2").Value = -0.01367558870849787
$ws.Range("AK2").Value = -0.02191235059760955
$ws.Range("AL2").Value = 36.2
$ws.Range("AM2").Value = 36.2
$ws.Range("AN2").Value = 0.7264472190692396
$ws.Range("AO2").Value = 4.516574585635359
$ws.Range("AP2").Value = -0.1061293984108966
$ws.Range("AQ2").Value = 4.516574585635359

$ws.Range("D3").Value = 0.121
$ws.Range("E3").Value = 0.07679999999999999
$ws.Range("G3").Value = 0.3078990592802747
$ws.Range("H3").Value = 0.3078990592802747
$ws.Range("I3").Value = 0.2441391667911005
$ws.Range("J3").Value = 0.155083590069419
$ws.Range("K3").Value = 98
$ws.Range("L3").Value = 0.1463341794833508
$ws.Range("M3").Value = 23.25
$ws.Range("N3").Value = 0.01677368155255754
$ws.Range("O3").Value = 0.2372448979591837
$ws.Range("P3").Value = 20.3
$ws.Range("Q3").Value = 0.01464540797922228
$ws.Range("R3").Value = 0.2071428571428572
$ws.Range("S3").Value = 2.949999999999999
$ws.Range("T3").Value = 0.1268817204301075
$ws.Range("U3").Value = 146.7
$ws.Range("V3").Value = 0.1058365197316211
$ws.Range("W3").Value = 0.1528622679769147
$ws.Range("X3").Value = 0.1129633422907675
$ws.Range("Y3").Value = 0.03989892568614717
$ws.Range("Z3").Value = 1.124223602484472
$ws.Range("AA3").Value = 0.1743486323140673
$ws.Range("AB3").Value = 0.108397913298502
$ws.Range("AC3").Value = 0.06595071901556529
$ws.Range("AD3").Value = 128
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 128
$ws.Range("AG3").Value = -18.69999999999999
$ws.Range("AH3").Value = 0.08453866983686679
$ws.Range("AI3").Value = 0.127987201279872
$ws.Range("AJ3").Value = -0.01367558870849787
$ws.Range("AK3").Value = -0.02191235059760955
$ws.Range("AL3").Value = 36.2
$ws.Range("AM3").Value = 36.2
$ws.Range("AN3").Value = 0.7264472190692396
$ws.Range("AO3").Value = 4.516574585635359
$ws.Range("AP3").Value = -0.1061293984108966
$ws.Range("AQ3").Value = 4.516574585635359

